$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 (L1) : update coordinates and rotation ---
$ws.Cells.Item(9, 4).Value  = "61.849mm"   # D9  Mid X
$ws.Cells.Item(9, 5).Value  = "-21.717mm"  # E9  Mid Y
$ws.Cells.Item(9, 6).Value  = "61.849mm"   # F9  Ref X
$ws.Cells.Item(9, 7).Value  = "-21.717mm"  # G9  Ref Y
$ws.Cells.Item(9, 8).Value  = "61.849mm"   # H9  Pad X
$ws.Cells.Item(9, 9).Value  = "-20.751mm"  # I9  Pad Y
$ws.Cells.Item(9, 12).Value = 270          # L9  Rotation

# --- Row 17 (U1) : update coordinates ---
$ws.Cells.Item(17, 4).Value = "58.42mm"    # D17 Mid X
$ws.Cells.Item(17, 5).Value = "-21.336mm"  # E17 Mid Y
$ws.Cells.Item(17, 6).Value = "58.42mm"    # F17 Ref X
$ws.Cells.Item(17, 7).Value = "-21.336mm"  # G17 Ref Y
$ws.Cells.Item(17, 8).Value = "59.691mm"   # H17 Pad X
$ws.Cells.Item(17, 9).Value = "-22.286mm"  # I17 Pad Y

# --- New Row 20 (EXP connector) ---
$ws.Cells.Item(20, 1).Value  = "EXP"
$ws.Cells.Item(20, 2).Value  = "SH1.0-6P"
$ws.Cells.Item(20, 3).Value  = "CONN-SMD_6P-P1.00-H-M_AFC10-S06QCC-00"
$ws.Cells.Item(20, 4).Value  = "64.643mm"
$ws.Cells.Item(20, 5).Value  = "-28.321mm"
$ws.Cells.Item(20, 6).Value  = "64.643mm"
$ws.Cells.Item(20, 7).Value  = "-28.321mm"
$ws.Cells.Item(20, 8).Value  = "62.143mm"
$ws.Cells.Item(20, 9).Value  = "-26.531mm"
$ws.Cells.Item(20, 10).Value = 8
$ws.Cells.Item(20, 11).Value = "T"
$ws.Cells.Item(20, 12).Value = 0
$ws.Cells.Item(20, 13).Value = "Yes"
$ws.Cells.Item(20, 14).Value = "SH1.0-6P"
